$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties in columns AC, AD, AE of row 1.
# Copy the formatting from the existing header cell (A1, style index 1 -
# bold font, border, centered/top alignment) onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record values for every player row (2-53): Wins=85, Losses=77, Ties=0
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 85  # AC
    $ws.Cells.Item($r, 30).Value = 77  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
